$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 12
$ws.Range("H12").Value = 542.8570999999999
$ws.Range("I12").Value = 766.6667
$ws.Range("J12").Value = 375
$ws.Range("K12").Value = 766.6667
$ws.Range("L12").Value = 375
$ws.Range("M12").Value = -596.6667
$ws.Range("N12").Value = -715

# ALC row 61
$ws.Range("H61").Value = 282.66666
$ws.Range("I61").Value = 282.66666
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 847.9999799999999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -675.9999799999999
$ws.Range("N61").ClearContents()

# ALC row 87
$ws.Range("H87").Value = 74816
$ws.Range("I87").Value = 74399
$ws.Range("J87").Value = 74899.39999999999
$ws.Range("K87").Value = 74399
$ws.Range("L87").Value = 74899.39999999999
$ws.Range("M87").Value = -73151

# ALC row 90
$ws.Range("H90").Value = 74816
$ws.Range("I90").Value = 74399
$ws.Range("J90").Value = 74899.39999999999
$ws.Range("K90").Value = 223197
$ws.Range("L90").Value = 224698.2
$ws.Range("M90").Value = -216957

# ALC row 98
$ws.Range("H98").Value = 1356.5
$ws.Range("I98").Value = 1356.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1356.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 141.5

# ALC row 100
$ws.Range("H100").Value = 2657.25
$ws.Range("I100").Value = 1252.5
$ws.Range("J100").Value = 3125.5
$ws.Range("K100").Value = 1252.5
$ws.Range("L100").Value = 3125.5
$ws.Range("M100").Value = -711.5

# ALC row 106
$ws.Range("H106").Value = 2406.3845
$ws.Range("I106").Value = 2116.818
$ws.Range("J106").Value = 3999
$ws.Range("K106").Value = 2116.818
$ws.Range("L106").Value = 3999
$ws.Range("M106").Value = -1485.818

# ALC row 122
$ws.Range("H122").Value = 1356.5
$ws.Range("I122").Value = 1356.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4069.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1619.5


$ws = $wb.Worksheets.Item("ARM")
# ARM row 5
$ws.Range("H5").Value = 243.16667
$ws.Range("I5").Value = 287.5
$ws.Range("J5").Value = 154.5
$ws.Range("K5").Value = 287.5
$ws.Range("L5").Value = 154.5
$ws.Range("M5").Value = -175.5

# ARM row 32
$ws.Range("H32").Value = 14918.875
$ws.Range("I32").Value = 7892.722
$ws.Range("J32").Value = 35997.332
$ws.Range("K32").Value = 7892.722
$ws.Range("L32").Value = 35997.332
$ws.Range("M32").Value = -7605.722

# ARM row 44
$ws.Range("H44").Value = 65776.60000000001
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 65776.60000000001
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 65776.60000000001
$ws.Range("N44").Value = -66752.60000000001

# ARM row 61
$ws.Range("H61").Value = 4857.1177
$ws.Range("I61").Value = 4065.5715
$ws.Range("J61").Value = 5411.2
$ws.Range("K61").Value = 4065.5715
$ws.Range("L61").Value = 5411.2
$ws.Range("M61").Value = -3853.5715
$ws.Range("N61").Value = -5835.2

# ARM row 74
$ws.Range("H74").Value = 5556.467
$ws.Range("I74").Value = 3620.7
$ws.Range("J74").Value = 9428
$ws.Range("K74").Value = 3620.7
$ws.Range("L74").Value = 9428
$ws.Range("M74").Value = -2746.7

# ARM row 77
$ws.Range("H77").Value = 5556.467
$ws.Range("I77").Value = 3620.7
$ws.Range("J77").Value = 9428
$ws.Range("K77").Value = 18103.5
$ws.Range("L77").Value = 47140
$ws.Range("M77").Value = -13735.5

# ARM row 136
$ws.Range("H136").Value = 4857.1177
$ws.Range("I136").Value = 4065.5715
$ws.Range("J136").Value = 5411.2
$ws.Range("K136").Value = 12196.7145
$ws.Range("L136").Value = 16233.6
$ws.Range("M136").Value = -9646.7145
$ws.Range("N136").Value = -21333.6


$ws = $wb.Worksheets.Item("BSM")
# BSM row 4
$ws.Range("H4").Value = 243.16667
$ws.Range("I4").Value = 287.5
$ws.Range("J4").Value = 154.5
$ws.Range("K4").Value = 287.5
$ws.Range("L4").Value = 154.5
$ws.Range("M4").Value = -172.5

# BSM row 60
$ws.Range("H60").Value = 19999.4
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 19999.4
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 19999.4
$ws.Range("N60").Value = -21197.4


$ws = $wb.Worksheets.Item("CRP")
# CRP row 9
$ws.Range("H9").Value = 68544.336
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 68544.336
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 68544.336
$ws.Range("N9").Value = -68880.336

# CRP row 31
$ws.Range("H31").Value = 5392.6924
$ws.Range("I31").Value = 3786.0476
$ws.Range("J31").Value = 7267.1113
$ws.Range("K31").Value = 3786.0476
$ws.Range("L31").Value = 7267.1113
$ws.Range("M31").Value = -3491.0476
$ws.Range("N31").Value = -7857.1113

# CRP row 34
$ws.Range("H34").Value = 5392.6924
$ws.Range("I34").Value = 3786.0476
$ws.Range("J34").Value = 7267.1113
$ws.Range("K34").Value = 3786.0476
$ws.Range("L34").Value = 7267.1113
$ws.Range("M34").Value = -3584.0476
$ws.Range("N34").Value = -7671.1113

# CRP row 38
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# CRP row 46
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# CRP row 99
$ws.Range("H99").Value = 4581.9
$ws.Range("I99").Value = 4356
$ws.Range("J99").Value = 4920.75
$ws.Range("K99").Value = 4356
$ws.Range("L99").Value = 4920.75
$ws.Range("M99").Value = -2858

# CRP row 122
$ws.Range("H122").Value = 4533.5557
$ws.Range("I122").Value = 3698
$ws.Range("J122").Value = 5202
$ws.Range("K122").Value = 11094
$ws.Range("L122").Value = 15606
$ws.Range("M122").Value = -8644
$ws.Range("N122").Value = -20506

# CRP row 126
$ws.Range("H126").Value = 4581.9
$ws.Range("I126").Value = 4356
$ws.Range("J126").Value = 4920.75
$ws.Range("K126").Value = 13068
$ws.Range("L126").Value = 14762.25
$ws.Range("M126").Value = -10598

# CRP row 132
$ws.Range("H132").Value = 2258.2222
$ws.Range("I132").Value = 1449.463
$ws.Range("J132").Value = 7110.778
$ws.Range("K132").Value = 4348.389
$ws.Range("L132").Value = 21332.334
$ws.Range("M132").Value = -1818.389

# CRP row 134
$ws.Range("H134").Value = 2017.1923
$ws.Range("I134").Value = 1671.3
$ws.Range("J134").Value = 3170.1667
$ws.Range("K134").Value = 5013.9
$ws.Range("L134").Value = 9510.500100000001
$ws.Range("M134").Value = -2478.9
$ws.Range("N134").Value = -14580.5001


$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 1287.3334
$ws.Range("I5").Value = 1194
$ws.Range("J5").Value = 1474
$ws.Range("K5").Value = 3582
$ws.Range("L5").Value = 4422
$ws.Range("M5").Value = -3470

# CUL row 107
$ws.Range("H107").Value = 434
$ws.Range("I107").Value = 330.75
$ws.Range("J107").Value = 485.625
$ws.Range("K107").Value = 992.25
$ws.Range("L107").Value = 1456.875
$ws.Range("M107").Value = 927.75
$ws.Range("N107").Value = -5296.875

# CUL row 113
$ws.Range("H113").Value = 943
$ws.Range("I113").Value = 516.6667
$ws.Range("J113").Value = 1198.8
$ws.Range("K113").Value = 1550.0001
$ws.Range("L113").Value = 3596.4
$ws.Range("M113").Value = 619.9999
$ws.Range("N113").Value = -7936.4

# CUL row 135
$ws.Range("H135").Value = 1287.3334
$ws.Range("I135").Value = 1194
$ws.Range("J135").Value = 1474
$ws.Range("K135").Value = 10746
$ws.Range("L135").Value = 13266
$ws.Range("M135").Value = -8211

# CUL row 137
$ws.Range("H137").Value = 5280.6875
$ws.Range("I137").Value = 2171.2222
$ws.Range("J137").Value = 9278.571
$ws.Range("K137").Value = 6513.6666
$ws.Range("L137").Value = 27835.713
$ws.Range("M137").Value = -1413.6666


$ws = $wb.Worksheets.Item("GSM")
# GSM row 131
$ws.Range("H131").Value = 70000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 70000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 70000
$ws.Range("N131").Value = -80080

# GSM row 134
$ws.Range("H134").Value = 58530
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 58530
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 175590
$ws.Range("N134").Value = -180660
$ws.Range("M134").ClearContents()


$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 3292.2856
$ws.Range("I22").Value = 4683
$ws.Range("J22").Value = 2249.25
$ws.Range("K22").Value = 4683
$ws.Range("L22").Value = 2249.25
$ws.Range("M22").Value = -4388
$ws.Range("N22").Value = -2839.25

# LTW row 27
$ws.Range("H27").Value = 3292.2856
$ws.Range("I27").Value = 4683
$ws.Range("J27").Value = 2249.25
$ws.Range("K27").Value = 4683
$ws.Range("L27").Value = 2249.25
$ws.Range("M27").Value = -4576
$ws.Range("N27").Value = -2463.25

# LTW row 40
$ws.Range("H40").Value = 11750.5
$ws.Range("I40").Value = 9376
$ws.Range("J40").Value = 16499.5
$ws.Range("K40").Value = 9376
$ws.Range("L40").Value = 16499.5
$ws.Range("M40").Value = -9240

# LTW row 46
$ws.Range("H46").Value = 1735.84
$ws.Range("I46").Value = 975.9
$ws.Range("J46").Value = 2242.4666
$ws.Range("K46").Value = 975.9
$ws.Range("L46").Value = 2242.4666
$ws.Range("M46").Value = -787.9
$ws.Range("N46").Value = -2618.4666

# LTW row 47
$ws.Range("H47").Value = 17500
$ws.Range("I47").Value = 10000
$ws.Range("J47").Value = 25000
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 25000
$ws.Range("M47").Value = -9510
$ws.Range("N47").Value = -25980

# LTW row 51
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# LTW row 52
$ws.Range("H52").Value = 17500
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 25000
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 25000
$ws.Range("M52").Value = -9767
$ws.Range("N52").Value = -25466

# LTW row 55
$ws.Range("H55").Value = 1404.8
$ws.Range("I55").Value = 424.25
$ws.Range("J55").Value = 2525.4285
$ws.Range("K55").Value = 424.25
$ws.Range("L55").Value = 2525.4285
$ws.Range("M55").Value = -251.25
$ws.Range("N55").Value = -2871.4285

# LTW row 122
$ws.Range("H122").Value = 6560.875
$ws.Range("I122").Value = 5581.1665
$ws.Range("J122").Value = 9500
$ws.Range("K122").Value = 16743.4995
$ws.Range("L122").Value = 28500
$ws.Range("M122").Value = -14293.4995


$ws = $wb.Worksheets.Item("WVR")
# WVR row 96
$ws.Range("H96").Value = 13800
$ws.Range("I96").Value = 4950
$ws.Range("J96").Value = 22650
$ws.Range("K96").Value = 4950
$ws.Range("L96").Value = 22650
$ws.Range("M96").Value = -3577
$ws.Range("N96").Value = -25396

# WVR row 107
$ws.Range("H107").Value = 1453.6285
$ws.Range("I107").Value = 1721.909
$ws.Range("J107").Value = 999.61536
$ws.Range("K107").Value = 5165.727000000001
$ws.Range("L107").Value = 2998.84608
$ws.Range("M107").Value = -3245.727000000001
$ws.Range("N107").Value = -6838.84608

# WVR row 122
$ws.Range("H122").Value = 4975.375
$ws.Range("I122").Value = 4300.6665
$ws.Range("J122").Value = 6999.5
$ws.Range("K122").Value = 12901.9995
$ws.Range("L122").Value = 20998.5
$ws.Range("M122").Value = -10451.9995

# WVR row 132
$ws.Range("H132").Value = 1907.5536
$ws.Range("I132").Value = 1537.3864
$ws.Range("J132").Value = 3264.8333
$ws.Range("K132").Value = 4612.1592
$ws.Range("L132").Value = 9794.499899999999
$ws.Range("M132").Value = -2082.1592

